# "add carolina as single"
#
# Re-shuffle the tail of the releases table (rows 102-105) and append a new
# row 106 for the "Carolina" single.
#
#   Before                                   After
#   102 B=Wildest Dreams   C=44456            102 B=Wildest Dreams   D=44456
#   103 B=Joker And Queen  C=44498            103 A=Red (TV)         D=44512
#   104 A=Red (TV)         D=44512            104 B=Joker And Queen  D=44603
#   105 B=This Love        C=44687            105 B=This Love        D=44687
#                                              106 B=Carolina         D=44736   (NEW)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateFmt = "yyyy\-mm\-dd;@"
$textFmt = "@"

# Row 102: keep "Wildest Dreams (Taylor's Version)" in B, move its release
# date from C102 into D102. Clear() (not ClearContents) so the vacated
# cell drops out of the XML entirely instead of leaving a bare <c s="1"/>.
$ws.Range("C102").Clear()
$ws.Range("B102").Value = "Wildest Dreams (Taylor's Version)"
$ws.Range("B102").NumberFormat = $textFmt
$ws.Range("D102").Value = 44456
$ws.Range("D102").NumberFormat = $dateFmt

# Row 103: now holds the "Red (Taylor's Version)" album entry that used to
# live in row 104 (album name in column A, date in column D).
$ws.Range("B103").Clear()
$ws.Range("C103").Clear()
$ws.Range("A103").Value = "Red (Taylor's Version)"
$ws.Range("D103").Value = 44512
$ws.Range("D103").NumberFormat = $dateFmt

# Row 104: now holds "The Joker And The Queen", with an updated release
# date (was C103=44498, now D104=44603), and no album name.
$ws.Range("A104").Clear()
$ws.Range("B104").Value = "The Joker And The Queen"
$ws.Range("B104").NumberFormat = $textFmt
$ws.Range("D104").Value = 44603
$ws.Range("D104").NumberFormat = $dateFmt

# Row 105: keep "This Love (Taylor's Version)" in B, move its release date
# from C105 into D105.
$ws.Range("C105").Clear()
$ws.Range("B105").Value = "This Love (Taylor's Version)"
$ws.Range("B105").NumberFormat = $textFmt
$ws.Range("D105").Value = 44687
$ws.Range("D105").NumberFormat = $dateFmt

# Row 106 (new): "Carolina" single.
$ws.Range("B106").Value = "Carolina"
$ws.Range("B106").NumberFormat = $textFmt
$ws.Range("D106").Value = 44736
$ws.Range("D106").NumberFormat = $dateFmt

# Match the final selection/active cell from the authored workbook.
$ws.Range("D106").Select()
